$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold/border/centered style for the A column label,
# plus the "blank cell" formatting used in columns E and V) from the
# existing last row (row 10) down into the new row 11, so the new row
# mirrors the existing layout before we fill in the 2021 figures.
$ws.Range("A10:AQ10").Copy()
$ws.Range("A11:AQ11").PasteSpecial(-4122)

# Now populate row 11 with the 2021 data (columns E11 and V11 stay
# blank, matching the source data for this row).
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 631.27
$ws.Range("C11").Value = 136.07
$ws.Range("D11").Value = 36.1
$ws.Range("F11").Value = 519.95
$ws.Range("G11").Value = 852.72
$ws.Range("H11").Value = 133.03
$ws.Range("I11").Value = 707.09
$ws.Range("J11").Value = 60.53
$ws.Range("K11").Value = 12194.63
$ws.Range("L11").Value = 93.88
$ws.Range("M11").Value = 8.56
$ws.Range("N11").Value = 0.1
$ws.Range("O11").Value = 515.04
$ws.Range("P11").Value = 175.9
$ws.Range("Q11").Value = 3.53
$ws.Range("R11").Value = 43.19
$ws.Range("S11").Value = 426.45
$ws.Range("T11").Value = 1.09
$ws.Range("U11").Value = 1377.08
$ws.Range("W11").Value = 25.43
$ws.Range("X11").Value = 57.45
$ws.Range("Y11").Value = 20.67
$ws.Range("Z11").Value = 894.63
$ws.Range("AA11").Value = 130.97
$ws.Range("AB11").Value = 85.13
$ws.Range("AC11").Value = 9.52
$ws.Range("AD11").Value = 276.68
$ws.Range("AE11").Value = 241.52
$ws.Range("AF11").Value = 2107.13
$ws.Range("AG11").Value = 900.9
$ws.Range("AH11").Value = 213.12
$ws.Range("AI11").Value = 131.46
$ws.Range("AJ11").Value = 9.4
$ws.Range("AK11").Value = 482.29
$ws.Range("AL11").Value = 113.67
$ws.Range("AM11").Value = 316.33
$ws.Range("AN11").Value = 3.86
$ws.Range("AO11").Value = 251.35
$ws.Range("AP11").Value = 193.95
$ws.Range("AQ11").Value = 7.1
